$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = "2024-06-15 08:15:34"
$ws.Range("D22").Value = 200
$ws.Range("E22").Value = 7

$ws.Range("A23").Value = 22
$ws.Range("B23").Value = 2
$ws.Range("C23").Value = "2024-06-15 08:15:34"
$ws.Range("D23").Value = 200
$ws.Range("E23").Value = 0
